$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New team member row: add repository hyperlink in D13 (Jose Alfredo Peto Martinez)
$ws.Range("D13").Value = "https://github.com/petoalfredo1/Tesis"
$ws.Hyperlinks.Add($ws.Range("D13"), "https://github.com/petoalfredo1/Tesis")
$ws.Range("D13").Style = "Hipervínculo"

# Row 18 had an oversized custom height left over from earlier formatting; restore
# it to the sheet's normal auto height.
$ws.Rows.Item(18).AutoFit()

# Move/save the active cell selection as it was left by the editor.
$ws.Range("I18").Select() | Out-Null
